$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.329.77"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").Value = "2.594.75"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'586.63"
$ws.Range("E5").Value = "  -3.06%  "

$ws.Range("D6").Value = "'148.95"
$ws.Range("E6").Value = "  -1.82%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'0.582"
$ws.Range("E8").Value = "  -1.61%  "

$ws.Range("D9").Value = "'0.108"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").Value = "'5.78"
$ws.Range("E10").Value = "  +1.79%  "

$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").Value = "'27.43"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("D14").Value = "3.060.51"
$ws.Range("E14").Value = "  -1.86%  "

$ws.Range("D15").Value = "63.167.29"
$ws.Range("E15").Value = "  -1.53%  "

$ws.Range("E16").Value = "  +3.41%  "

$ws.Range("D17").Value = "2.577.03"
$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("D18").Value = "'11.99"
$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("D19").Value = "'4.64"
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").Value = "'343.03"
$ws.Range("E20").Value = "  -2.85%  "

$ws.Range("D21").Value = "'6.80"
$ws.Range("E21").Value = "  -2.41%  "

$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").Value = "'66.36"

$ws.Range("D24").Value = "'1.71"
$ws.Range("E24").Value = "  -2.32%  "

$ws.Range("D25").Value = "'9.16"
$ws.Range("E25").Value = "  -2.74%  "

$ws.Range("E26").Value = "  -5.13%  "

$ws.Range("D27").Value = "'565.80"
$ws.Range("E27").Value = "  +3.65%  "

$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.161"
$ws.Range("E30").Value = "  -3.53%  "

$ws.Range("E31").Value = "  -2.67%  "

$ws.Range("D32").Value = "0.0₃0837"
$ws.Range("E32").Value = "  -3.26%  "

$ws.Range("D33").Value = "'1.75"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").Value = "'5.27"
$ws.Range("E34").Value = "  -0.94%  "

$ws.Range("D35").Value = "'165.44"
$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").Value = "'19.31"
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("D39").Value = "'1.91"
$ws.Range("E39").Value = "  -5.52%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").Value = "'165.66"
$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("D43").Value = "'22.65"
$ws.Range("E43").Value = "  +4.17%  "

$ws.Range("D44").Value = "'0.0576"
$ws.Range("E44").Value = "  -1.57%  "

$ws.Range("D45").Value = "'2.09"
$ws.Range("E45").Value = "  +1.87%  "

$ws.Range("D46").Value = "'0.628"

$ws.Range("D47").Value = "'0.0245"
$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("D48").Value = "'0.0955"
$ws.Range("E48").Value = "  -1.33%  "

$ws.Range("D49").Value = "'18.94"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("D50").Value = "0.0₆0225"
$ws.Range("E50").Value = "  +13.29%  "

$ws.Range("E51").Value = "  -4.34%  "
